$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target dataset for rows 2-17 (A:H)
# Columns: A=index, B=Chr, C=Position, D=refDepth, E=altDepth, F=sampleNames, G=VAF, H=chr_pos
$data = @(
    @(0,  "chr10", 500,  43, 7,  "Pt0_blast",  0.14, "chr10_500"),
    @(1,  "chr9",  127,  43, 7,  "Pt0_blast",  0.14, "chr9_127"),
    @(2,  "chr4",  502,  47, 3,  "Pt0_blast",  0.06, "chr4_502"),
    @(3,  "chrY",  1043, 50, 0,  "Pt0_blast",  0,    "chrY_1043"),
    @(4,  "chr10", 500,  34, 16, "Pt0_Er1",    0.32, "chr10_500"),
    @(5,  "chr9",  127,  48, 2,  "Pt0_Er1",    0.04, "chr9_127"),
    @(6,  "chr4",  502,  50, 0,  "Pt0_Er1",    0,    "chr4_502"),
    @(7,  "chrY",  1043, 50, 0,  "Pt0_Er1",    0,    "chrY_1043"),
    @(8,  "chr10", 500,  45, 5,  "Pt0_NK",     0.1,  "chr10_500"),
    @(9,  "chr9",  127,  49, 1,  "Pt0_NK",     0.02, "chr9_127"),
    @(10, "chr4",  502,  44, 6,  "Pt0_NK",     0.12, "chr4_502"),
    @(11, "chrY",  1043, 50, 0,  "Pt0_NK",     0,    "chrY_1043"),
    @(12, "chr10", 500,  50, 0,  "Pt0_Bcells", 0,    "chr10_500"),
    @(13, "chr9",  127,  50, 0,  "Pt0_Bcells", 0,    "chr9_127"),
    @(14, "chr4",  502,  50, 0,  "Pt0_Bcells", 0,    "chr4_502"),
    @(15, "chrY",  1043, 41, 9,  "Pt0_Bcells", 0.18, "chrY_1043")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
}

# Apply the existing index-column cell format (bold/centered/bordered, style id 1)
# to the newly added rows 11-17 so they match rows 2-10.
$ws.Range("A2").Copy()
$ws.Range("A11:A17").PasteSpecial(-4122)
$excel.CutCopyMode = 0
